$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Addr, $Val)
    $r = $ws.Range($Addr)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "67.070.54"

$ws.Range("D3").Value = "3.453.12"
$ws.Range("E3").Value = "  -1.72%  "

$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue "D5" "592.61"
$ws.Range("E5").Value = "  -1.36%  "

Set-TextValue "D6" "179.12"
$ws.Range("E6").Value = "  +1.48%  "

Set-TextValue "D7" "0.608"
$ws.Range("E7").Value = "  +3.07%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "3.449.97"
$ws.Range("E9").Value = "  -1.69%  "

Set-TextValue "D10" "0.138"
$ws.Range("E10").Value = "  +4.99%  "

$ws.Range("E11").Value = "  -2.75%  "

$ws.Range("E12").Value = "  -1.06%  "

$ws.Range("D13").Value = "4.054.41"
$ws.Range("E13").Value = "  -1.60%  "

Set-TextValue "D14" "31.93"
$ws.Range("E14").Value = "  +1.62%  "

$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("D16").Value = "67.062.14"
$ws.Range("E16").Value = "  -0.60%  "

Set-TextValue "D17" "0.0000176"
$ws.Range("E17").Value = "  -1.99%  "

$ws.Range("D18").Value = "3.458.06"
$ws.Range("E18").Value = "  -1.33%  "

Set-TextValue "D19" "6.18"
$ws.Range("E19").Value = "  -2.22%  "

Set-TextValue "D20" "14.10"
$ws.Range("E20").Value = "  -3.99%  "

Set-TextValue "D21" "391.39"
$ws.Range("E21").Value = "  -0.91%  "

Set-TextValue "D22" "7.89"
$ws.Range("E22").Value = "  -2.08%  "

$ws.Range("E23").Value = "  +1.10%  "

Set-TextValue "D24" "0.997"
$ws.Range("E24").Value = "  -0.12%  "

Set-TextValue "D25" "0.536"
$ws.Range("E25").Value = "  -0.91%  "

Set-TextValue "D26" "71.51"
$ws.Range("E26").Value = "  -2.81%  "

$ws.Range("E27").Value = "  -2.06%  "

Set-TextValue "D28" "10.32"
$ws.Range("E28").Value = "  +0.21%  "

$ws.Range("E29").Value = "  -3.45%  "

Set-TextValue "D30" "1.00"
$ws.Range("E30").Value = "  +0.54%  "

Set-TextValue "D31" "6.10"
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("E32").Value = "  -1.40%  "

$ws.Range("E33").Value = "  -2.91%  "

Set-TextValue "D34" "23.42"
$ws.Range("E34").Value = "  -1.50%  "

Set-TextValue "D35" "7.29"
$ws.Range("E35").Value = "  -1.54%  "

$ws.Range("E36").Value = "  -0.06%  "

Set-TextValue "D37" "1.57"
$ws.Range("E37").Value = "  -4.49%  "

Set-TextValue "D38" "160.84"
$ws.Range("E38").Value = "  -1.80%  "

Set-TextValue "D39" "0.874"
$ws.Range("E39").Value = "  -0.47%  "

$ws.Range("E40").Value = "  +9.36%  "

$ws.Range("E41").Value = "  -4.50%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D42" "6.73"
$ws.Range("E42").Value = "  -4.56%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D43" "4.64"
$ws.Range("E43").Value = "  -1.26%  "

Set-TextValue "D44" "25.98"
$ws.Range("E44").Value = "  -2.38%  "

Set-TextValue "D45" "0.0717"
$ws.Range("E45").Value = "  -2.50%  "

$ws.Range("D46").Value = "2.751.90"
$ws.Range("E46").Value = "  -1.98%  "

Set-TextValue "D47" "26.19"
$ws.Range("E47").Value = "  -4.13%  "

Set-TextValue "D48" "41.30"
$ws.Range("E48").Value = "  -2.99%  "

Set-TextValue "D49" "0.0297"
$ws.Range("E49").Value = "  -1.54%  "

Set-TextValue "D50" "323.40"
$ws.Range("E50").Value = "  -4.59%  "

$ws.Range("E51").Value = "  -4.27%  "
